# Populate the results sheet with the full scraped job listing (rows 2-16),
# extending the single sample row that was already present (row 2's A/B/C/D/E
# values get overwritten with the real first record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
  "new2022 Summer Data Analyst Intern",
  "newFX Data Analytics Internship, Remote, Summer/Fall 2022",
  "newData Center Platform Application Engineer Intern",
  "new2022 Summer PHD Data Analytics Research Intern",
  "newDeep Learning Software Engineer Graduate Internship",
  "newIOTG Research AI Scientist internship",
  "newSummer 2022 Data Science Intern",
  "newData Science Internship (Philadelphia) - Publicis Health",
  "newInventory Planning & Business Analytics Intern",
  "newData Science Summer Intern - AI Innovations",
  "newData Analytics Intern (Summer 2022)",
  "newResearch/Data Scientist Intern",
  "newData Science Intern",
  "newJunior Business Analyst Intern",
  "newIntern: Energy Trading Analyst"
)

$companies = @(
  "General Motors",
  "General Entertainment Content",
  "Intel",
  "General Motors",
  "Intel",
  "Intel",
  "Slack",
  "Publicis Health",
  "Fullbeauty",
  "IBM",
  "Poshmark",
  "Ascension",
  "Varian Medical Systems",
  "Elsevier",
  "Greenwich Commodities LLC"
)

$links = @(
  "www.indeed.com//cmp/General-Motors",
  "www.indeed.com//q-General-Entertainment-Content-l-Burbank,-CA-jobs.html",
  "www.indeed.com//cmp/Intel-Corporation",
  "www.indeed.com//cmp/General-Motors",
  "www.indeed.com//cmp/Intel-Corporation",
  "www.indeed.com//cmp/Intel-Corporation",
  "www.indeed.com//cmp/Slack",
  "www.indeed.com//cmp/Publicis-Healthcare-Communications-Group",
  "www.indeed.com//cmp/Fullbeauty",
  "www.indeed.com//cmp/IBM",
  "www.indeed.com//cmp/Poshmark",
  "www.indeed.com//cmp/Ascension",
  "www.indeed.com//cmp/Varian-Medical-Systems",
  "www.indeed.com//cmp/Relx-Group",
  "www.indeed.com//jobs?q=Greenwich+Commodities+LLC&l=Denver,+CO&nc=jasx"
)

$dates = @(
  "PostedToday","PostedToday","PostedToday","PostedToday","PostedToday",
  "PostedToday","PostedToday","PostedToday","PostedToday","PostedToday",
  "PostedToday","PostedToday","PostedToday","PostedToday","PostedToday"
)

$rowCount = $titles.Length

# Give every data row (A2:A16) the same style as the original sample row
# so the numbering column keeps its formatting (bordered/centered style).
for ($i = 1; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Range("A2").Copy($ws.Range("A$r"))
}

# Column B: titles (write the whole column first - matches scrape order)
for ($i = 0; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 2).Value = $titles[$i]
}

# Column C: companies
for ($i = 0; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 3).Value = $companies[$i]
}

# Column D: links
for ($i = 0; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 4).Value = $links[$i]
}

# Column E: date_listed
for ($i = 0; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 5).Value = $dates[$i]
}

# Row index numbers in column A: row2 stays 0 (already there), rows 3-16 get 1..14
for ($i = 1; $i -lt $rowCount; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 1).Value = $i
}
